$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.624.58"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "2.489.35"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.69"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.05"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "2.512.16"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").Value = "2.936.84"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.04"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "58.589.67"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "2.504.42"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.22"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.44"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.53"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0767"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.62"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -2.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "167.33"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.14"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.45"
$ws.Range("E35").Value = "  -0.29%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.37"
$ws.Range("E36").Value = "  -5.52%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -2.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.67"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.812"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.03"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.85"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.29"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0510"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.24"
$ws.Range("E51").Value = "  -1.85%  "
